# Append the new daily mod-count entry as row 84.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 83
$newRow = 84

# Copy the formatting of the preceding data row so the new row matches
# (centered alignment style used throughout the data rows).
$ws.Range("A$lastRow`:C$lastRow").Copy()
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122) # xlPasteFormats

# Column A holds the date as plain text (e.g. "2026/01/31"), not a real
# date value, so force Text formatting before assigning it to avoid Excel
# auto-converting the string into a date serial number.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/02/02"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1170

# Re-apply the shared formatting to column A so it matches the sibling
# cells exactly (the Text-format step above nudges its style off of theirs).
$ws.Range("B$newRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122) # xlPasteFormats
